$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pairs of rows whose data (columns B:AD) must be swapped with each other.
# Column A (the running row index) stays untouched.
$pairs = @(
    @(78, 79),
    @(98, 99),
    @(134, 135),
    @(187, 188),
    @(196, 197),
    @(365, 366),
    @(367, 368)
)

foreach ($pair in $pairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    $rng1 = $ws.Range("B$r1`:AD$r1")
    $rng2 = $ws.Range("B$r2`:AD$r2")

    $vals1 = $rng1.Value2
    $vals2 = $rng2.Value2

    $rng1.Value = $vals2
    $rng2.Value = $vals1
}
